$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.005464
$ws.Range("H2").Value = 0.016392
$ws.Range("I2").Value = 0.002556359763658365
$ws.Range("J2").Value = 0.002556359763658364
$ws.Range("M2").Value = 0.097952
$ws.Range("N2").Value = 0.293856
$ws.Range("O2").Value = 0.3056806443660103
$ws.Range("P2").Value = 0.3056806443660104
$ws.Range("Q2").Value = 0.0005352097280000001
$ws.Range("R2").Value = 0.004816887552
$ws.Range("S2").Value = 0.0007814296997864308
$ws.Range("T2").Value = 0.0007814296997864308
$ws.Range("G3").Value = 0.005464
$ws.Range("H3").Value = 0.016392
$ws.Range("I3").Value = 0.002556359763658365
$ws.Range("J3").Value = 0.002556359763658364
$ws.Range("O3").Value = 0.3650866467564809
$ws.Range("P3").Value = 0.3650866467564809
$ws.Range("Q3").Value = 0.000639222432
$ws.Range("R3").Value = 0.005753001888
$ws.Range("S3").Value = 0.0009332928140172224
$ws.Range("T3").Value = 0.0009332928140172223
$ws.Range("G4").Value = 0.005464
$ws.Range("H4").Value = 0.016392
$ws.Range("I4").Value = 0.002556359763658365
$ws.Range("J4").Value = 0.002556359763658364
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.105499
$ws.Range("N4").Value = 0.316497
$ws.Range("O4").Value = 0.3292327088775087
$ws.Range("P4").Value = 0.3292327088775087
$ws.Range("Q4").Value = 0.000576446536
$ws.Range("R4").Value = 0.005188018824
$ws.Range("S4").Value = 0.0008416372498547113
$ws.Range("T4").Value = 0.0008416372498547113
$ws.Range("I5").Value = 0.8402398349532294
$ws.Range("J5").Value = 0.8402398349532293
$ws.Range("M5").Value = 0.097952
$ws.Range("N5").Value = 0.293856
$ws.Range("O5").Value = 0.3056806443660103
$ws.Range("P5").Value = 0.3056806443660104
$ws.Range("Q5").Value = 0.1759159801813333
$ws.Range("R5").Value = 1.583243821632
$ws.Range("S5").Value = 0.2568450541704934
$ws.Range("T5").Value = 0.2568450541704934
$ws.Range("I6").Value = 0.8402398349532294
$ws.Range("J6").Value = 0.8402398349532293
$ws.Range("O6").Value = 0.3650866467564809
$ws.Range("P6").Value = 0.3650866467564809
$ws.Range("S6").Value = 0.3067603438142935
$ws.Range("T6").Value = 0.3067603438142935
$ws.Range("I7").Value = 0.8402398349532294
$ws.Range("J7").Value = 0.8402398349532293
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.105499
$ws.Range("N7").Value = 0.316497
$ws.Range("O7").Value = 0.3292327088775087
$ws.Range("P7").Value = 0.3292327088775087
$ws.Range("Q7").Value = 0.1894699443926667
$ws.Range("R7").Value = 1.705229499534
$ws.Range("S7").Value = 0.2766344369684425
$ws.Range("T7").Value = 0.2766344369684425
$ws.Range("G8").Value = 0.3360096666666667
$ws.Range("H8").Value = 1.008029
$ws.Range("I8").Value = 0.1572038052831124
$ws.Range("J8").Value = 0.1572038052831123
$ws.Range("M8").Value = 0.097952
$ws.Range("N8").Value = 0.293856
$ws.Range("O8").Value = 0.3056806443660103
$ws.Range("P8").Value = 0.3056806443660104
$ws.Range("Q8").Value = 0.03291281886933334
$ws.Range("R8").Value = 0.296215369824
$ws.Range("S8").Value = 0.04805416049573061
$ws.Range("T8").Value = 0.04805416049573061
$ws.Range("G9").Value = 0.3360096666666667
$ws.Range("H9").Value = 1.008029
$ws.Range("I9").Value = 0.1572038052831124
$ws.Range("J9").Value = 0.1572038052831123
$ws.Range("O9").Value = 0.3650866467564809
$ws.Range("P9").Value = 0.3650866467564809
$ws.Range("Q9").Value = 0.03930909888400001
$ws.Range("R9").Value = 0.353781889956
$ws.Range("S9").Value = 0.05739301012817025
$ws.Range("T9").Value = 0.05739301012817025
$ws.Range("G10").Value = 0.3360096666666667
$ws.Range("H10").Value = 1.008029
$ws.Range("I10").Value = 0.1572038052831124
$ws.Range("J10").Value = 0.1572038052831123
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.105499
$ws.Range("N10").Value = 0.316497
$ws.Range("O10").Value = 0.3292327088775087
$ws.Range("P10").Value = 0.3292327088775087
$ws.Range("Q10").Value = 0.03544868382366667
$ws.Range("R10").Value = 0.319038154413
$ws.Range("S10").Value = 0.05175663465921149
$ws.Range("T10").Value = 0.05175663465921149

Write-Output "Applied 100 cell updates"